$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulated-game transition probabilities for Sheet1 (New Mexico St._A matrix)
# after adding more games and refining the simulation logic.

$ws.Range("B2").Value = 0.2244897959183673
$ws.Range("C2").Value = 0.5340136054421769
$ws.Range("J2").Value = 0.006802721088435374
$ws.Range("P2").Value = 0.1326530612244898
$ws.Range("S2").Value = 0.1020408163265306
$ws.Range("B3").Value = 0.006172839506172839
$ws.Range("C3").Value = 0.0308641975308642
$ws.Range("J3").Value = 0.0308641975308642
$ws.Range("P3").Value = 0.7160493827160493
$ws.Range("S3").Value = 0.2160493827160494
$ws.Range("J4").Value = 0.02439024390243903
$ws.Range("P4").Value = 0.6585365853658537
$ws.Range("S4").Value = 0.3170731707317073
$ws.Range("B6").Value = 0.07798165137614679
$ws.Range("D6").Value = 0.01376146788990826
$ws.Range("E6").Value = 0.004587155963302753
$ws.Range("F6").Value = 0.06422018348623854
$ws.Range("J6").Value = 0.2155963302752294
$ws.Range("O6").Value = 0.03211009174311927
$ws.Range("Q6").Value = 0.1330275229357798
$ws.Range("R6").Value = 0.06880733944954129
$ws.Range("S6").Value = 0.3899082568807339
$ws.Range("B7").Value = 0.1106382978723404
$ws.Range("D7").Value = 0.02553191489361702
$ws.Range("F7").Value = 0.0425531914893617
$ws.Range("J7").Value = 0.1574468085106383
$ws.Range("O7").Value = 0.02127659574468085
$ws.Range("Q7").Value = 0.2212765957446808
$ws.Range("R7").Value = 0.08085106382978724
$ws.Range("S7").Value = 0.3404255319148936
$ws.Range("B8").Value = 0.1061946902654867
$ws.Range("D8").Value = 0.01474926253687316
$ws.Range("E8").Value = 0.002949852507374631
$ws.Range("F8").Value = 0.05604719764011799
$ws.Range("J8").Value = 0.08554572271386431
$ws.Range("O8").Value = 0.02949852507374631
$ws.Range("Q8").Value = 0.191740412979351
$ws.Range("R8").Value = 0.06784660766961652
$ws.Range("S8").Value = 0.4454277286135693
$ws.Range("B9").Value = 0.07731958762886598
$ws.Range("D9").Value = 0.01030927835051546
$ws.Range("F9").Value = 0.06701030927835051
$ws.Range("J9").Value = 0.1288659793814433
$ws.Range("O9").Value = 0.04123711340206185
$ws.Range("Q9").Value = 0.1907216494845361
$ws.Range("R9").Value = 0.05154639175257732
$ws.Range("S9").Value = 0.4329896907216495
$ws.Range("B10").Value = 0.1093355761143818
$ws.Range("D10").Value = 0.02523128679562658
$ws.Range("F10").Value = 0.07401177460050462
$ws.Range("J10").Value = 0.1051303616484441
$ws.Range("O10").Value = 0.02775441547518924
$ws.Range("Q10").Value = 0.1925988225399495
$ws.Range("R10").Value = 0.07821698906644239
$ws.Range("S10").Value = 0.3877207737594617
$ws.Range("F11").Value = 0.002695417789757413
$ws.Range("G11").Value = 0.1536388140161725
$ws.Range("J11").Value = 0.08355795148247978
$ws.Range("K11").Value = 0.2318059299191375
$ws.Range("L11").Value = 0.5148247978436657
$ws.Range("S11").Value = 0.01347708894878706
$ws.Range("G12").Value = 0.7058823529411765
$ws.Range("J12").Value = 0.2254901960784314
$ws.Range("K12").Value = 0.01470588235294118
$ws.Range("L12").Value = 0.02450980392156863
$ws.Range("G13").Value = 0.7843137254901961
$ws.Range("J13").Value = 0.1568627450980392
$ws.Range("S13").Value = 0.05882352941176471
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.04460966542750929
$ws.Range("H15").Value = 0.1263940520446097
$ws.Range("I15").Value = 0.09665427509293681
$ws.Range("J15").Value = 0.3011152416356878
$ws.Range("K15").Value = 0.104089219330855
$ws.Range("M15").Value = 0.01486988847583643
$ws.Range("O15").Value = 0.104089219330855
$ws.Range("S15").Value = 0.20817843866171
$ws.Range("F16").Value = 0.02272727272727273
$ws.Range("H16").Value = 0.125
$ws.Range("I16").Value = 0.09659090909090909
$ws.Range("J16").Value = 0.4090909090909091
$ws.Range("K16").Value = 0.1420454545454546
$ws.Range("M16").Value = 0.01136363636363636
$ws.Range("O16").Value = 0.08522727272727272
$ws.Range("S16").Value = 0.1079545454545455
$ws.Range("F17").Value = 0.01477832512315271
$ws.Range("H17").Value = 0.1428571428571428
$ws.Range("I17").Value = 0.07881773399014778
$ws.Range("J17").Value = 0.4433497536945813
$ws.Range("K17").Value = 0.1157635467980296
$ws.Range("M17").Value = 0.01970443349753695
$ws.Range("O17").Value = 0.05911330049261083
$ws.Range("S17").Value = 0.125615763546798
$ws.Range("F18").Value = 0.02484472049689441
$ws.Range("H18").Value = 0.1490683229813665
$ws.Range("I18").Value = 0.1055900621118012
$ws.Range("J18").Value = 0.4285714285714285
$ws.Range("K18").Value = 0.124223602484472
$ws.Range("M18").Value = 0.0124223602484472
$ws.Range("N18").Value = 0.006211180124223602
$ws.Range("O18").Value = 0.08695652173913043
$ws.Range("S18").Value = 0.06211180124223602
$ws.Range("F19").Value = 0.02289452166802943
$ws.Range("H19").Value = 0.169255928045789
$ws.Range("I19").Value = 0.08503679476696648
$ws.Range("J19").Value = 0.3638593622240393
$ws.Range("K19").Value = 0.1300081766148814
$ws.Range("M19").Value = 0.02861815208503679
$ws.Range("N19").Value = 0.001635322976287817
$ws.Range("O19").Value = 0.0776778413736713
$ws.Range("S19").Value = 0.1210139002452984
